# Resize/reposition the timeline text placeholders and bump the font
# sizes of their text so the enlarged boxes read well (timeline resize).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- "Ransundet grundas av ..." (ph idx=14) ---------------------------
$sh = $s.Shapes.Item(14)
$sh.Left   = 34.617559455118105
$sh.Top    = 361.8467865535433
$sh.Width  = 169.99992125984252
$sh.Height = 47.37496062992126
$tr = $sh.TextFrame.TextRange
$tr.Runs(1).Font.Size = 18
$tr.Runs(2).Font.Size = 18

# --- "Vag bryts till byn ..." (ph idx=15) ------------------------------
$sh = $s.Shapes.Item(15)
$sh.Left   = 281.90055848110234
$sh.Top    = 361.8467865535433
$sh.Width  = 156.19889763779528
$sh.Height = 47.37496062992126
$sh.TextFrame.TextRange.Runs(1).Font.Size = 18

# --- "Ransundets vilthagn invigs ..." (ph idx=16) ----------------------
$sh = $s.Shapes.Item(16)
$sh.Left   = 479.8087401574803
$sh.Top    = 361.63929753858264
$sh.Width  = 239.99992125984252
$sh.Height = 47.37496062992126
$sh.TextFrame.TextRange.Runs(1).Font.Size = 18

# --- "Gatubelysning satts upp i byn." (ph idx=17) ----------------------
$sh = $s.Shapes.Item(17)
$sh.Left   = 761.9004212007874
$sh.Top    = 361.63929753858264
$sh.Width  = 156.19889763779528
$sh.Height = 47.37496062992126
$sh.TextFrame.TextRange.Runs(1).Font.Size = 18

# --- Year labels: bump font size to 32pt (ph idx=18..21) ---------------
$s.Shapes.Item(18).TextFrame.TextRange.Runs(1).Font.Size = 32
$s.Shapes.Item(19).TextFrame.TextRange.Runs(1).Font.Size = 32
$s.Shapes.Item(20).TextFrame.TextRange.Runs(1).Font.Size = 32
$s.Shapes.Item(21).TextFrame.TextRange.Runs(1).Font.Size = 32
